# Daily attendance processing - swap the order of names in the
# "Recorded By" (column G) cells that read "dnasr281@gmail.com, System"
# so they read "System, dnasr281@gmail.com" instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = "dnasr281@gmail.com, System"
$replacement = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $target) {
        $cell.Value = $replacement
    }
}
